# "Add files via upload" — refresh of data/params_scenarios.xlsx.
# On the scenario_params_paper sheet the probability parameters in column O
# (was 0.6) and column P (was 0.3) are lowered to 0.5 / 0.2 respectively.
# Row 2 holds the literal values; rows 3-17 hold shared formulas
# ("+O$2" / "+P$2") that recompute automatically once row 2 changes.
# The previously hidden helper columns H:AA are also unhidden, and the
# sheet's stored selection moves to P3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_params_paper")

$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 0.2

$ws.Columns("H:AA").Hidden = $false

$ws.Activate()
$ws.Range("P3").Select()
